# Commit: "corrections for new lines in xls and tests"
#
# Sheet1 gets a new blank row inserted above the old row 3 (the row that held
# the CONCAT/sum/literal formulas). That pushes the old row 3 down to row 4
# and grows the sheet's dimension from A1:C3 to A1:C4. The freshly inserted
# row 3 picks up the formatting of the row above it (style index 1, same as
# row 2), matching Excel's default "inherit formatting from the row above"
# behaviour on a row insert.
#
# The selections on both sheets are also updated to reflect the row-3
# insertion point that was used to drive the edit.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")

# Insert a new row above row 3 on Sheet1 -- old row 3 (formulas) becomes row 4.
$ws1.Rows.Item(3).EntireRow.Insert()

# Reflect the row-3 insertion point in the saved selection state. Sheet2's
# data is untouched, but its selection now spans the same row-3 plus the
# pre-existing E15 cell (active cell stays E15).
$ws2.Activate()
$ws2.Range("E15,3:3").Select()

$ws1.Activate()
$ws1.Range("A3").EntireRow.Select()
